# Update "想去人数" (want-to-go count) figures that changed between the
# two data pulls recorded in the workbook.
#
# Sheet "展览" and sheet "全部类型" both contain the same rows of
# exhibition data in column F (rows 2, 3, 10, 11, 16). Update both.

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 1141
    3  = 583
    10 = 5171
    11 = 4767
    16 = 183
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
